$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.807.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.491.07'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '534.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.61%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.561'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.511.80'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.101'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.41%  '
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('E12').Value = '  -2.41%  '
$ws.Range('E13').Value = '  -1.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.938.53'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.06'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.802.52'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.507.82'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.00'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.67'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.52%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.50'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.74%  '
$ws.Range('E25').Value = '  -0.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.164'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.995'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.56'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.72'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.80%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('E35').Value = '  -3.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.45'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('E37').Value = '  -2.80%  '
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.67'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.817'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.62'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '277.64'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.995'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.84'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '126.11'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0923'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0510'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('E50').Value = '  -1.53%  '
$ws.Range('E51').Value = '  -0.65%  '
